# forests-scraped.xlsx update — 2025-10-27 12:18
#
# The two listings that were sitting in the "New" sheet get archived into
# "Previously added" (appended at the bottom), and a fresh batch of 9
# scraped listings replaces the "New" sheet's contents.

$wb = $excel.ActiveWorkbook
$wsOld = $wb.Worksheets.Item("Previously added")
$wsNew = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# 1. Archive the current "New" rows (2 & 3) onto the end of
#    "Previously added" (rows 208 & 209), keeping their hyperlinks.
# ---------------------------------------------------------------------

$archiveRows = @(
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/maltas-pag/gipgg.html"; Price="21 600 €";  District="Rēzekne un raj."; Area="5.40 ha."; Cadastre="78700070207"; Date=45954.56736111111 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/tukums-and-reg/ceres-pag/ffcml.html";    Price="70 €";      District="Tukums un raj.";  Area="2 ha.";    Cadastre="90440010002"; Date=45953.91736111111 }
)

# Seed rows 208:209 with the same look as the last existing data row
# (207) before typing anything into them.
$wsOld.Range("A207:F207").Copy()
$wsOld.Range("A208:F209").PasteSpecial(-4122)   # xlPasteFormats
$wsOld.Application.CutCopyMode = $false

# Force plain-text entry for the price/cadastre columns so digit-only
# values such as "70 €" or "90440010002" aren't auto-converted to
# numbers on write.
$wsOld.Range("B208:B209").NumberFormat = "@"
$wsOld.Range("E208:E209").NumberFormat = "@"

$destRow = 208
foreach ($row in $archiveRows) {
    $wsOld.Range("A" + $destRow).Value = $row.Link
    $wsOld.Range("B" + $destRow).Value = $row.Price
    $wsOld.Range("C" + $destRow).Value = $row.District
    $wsOld.Range("D" + $destRow).Value = $row.Area
    $wsOld.Range("E" + $destRow).Value = $row.Cadastre
    $wsOld.Range("F" + $destRow).Value = $row.Date

    $wsOld.Hyperlinks.Add($wsOld.Range("A" + $destRow), $row.Link)

    $destRow = $destRow + 1
}

# Hyperlinks.Add stamps column A with Excel's built-in "Hyperlink" look;
# restore the sheet's own underline style used by every other link cell.
$wsOld.Range("A207").Copy()
$wsOld.Range("A208:A209").PasteSpecial(-4122)   # xlPasteFormats
$wsOld.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Replace the "New" sheet's contents with the freshly scraped batch.
# ---------------------------------------------------------------------

# Drop the two old hyperlinks so they don't linger on cells whose text
# is about to change.
$wsNew.Range("A2").Hyperlinks.Delete()
$wsNew.Range("A3").Hyperlinks.Delete()

$newRows = @(
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/ziguru-pag/dfxge.html";       Price="9 000 €";   District="Balvi un raj.";     Area="1 ha.";      Cadastre="38980010052"; Date=45957.40625 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/raunas-pag/jjxpi.html";        Price="25 000 €";  District="Cēsis un raj.";     Area="3 ha.";      Cadastre="42760080028"; Date=45957.49444444444 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/garsenes-pag/bxojlf.html"; Price="57 000 €";  District="Jēkabpils un raj."; Area="13 ha.";     Cadastre="56620020013"; Date=45955.50069444445 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kalniesu-pag/jxgdk.html";   Price="50 000 €";  District="Krāslava un raj.";  Area="10 ha.";     Cadastre="60680040608"; Date=45955.45 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/liepaja-and-reg/durbe/gdkpm.html";           Price="125 000 €"; District="Liepāja un raj.";   Area="35 ha.";     Cadastre="64270060039"; Date=45955.42152777778 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/livani/jlngb.html";           Price="22 500 €";  District="Preiļi un raj.";    Area="4 ha.";      Cadastre="76860070184"; Date=45957.44236111111 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/kaunatas-pag/npjbg.html";    Price="13 000 €";  District="Rēzekne un raj.";   Area="4 ha.";      Cadastre="78620090056"; Date=45957.57430555555 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/feimanu-pag/mhedc.html";     Price="50 000 €";  District="Rēzekne un raj.";   Area="11.83 ha.";  Cadastre="78520030192"; Date=45957.46805555555 },
    @{ Link="https://www.ss.com/msg/lv/real-estate/wood/talsi-and-reg/kulciema-pag/cdlld.html";      Price="51 000 €";  District="Talsi un raj.";     Area="12 ha.";     Cadastre="88640040079"; Date=45957.40833333333 }
)

# Seed rows 2:10 with the same look as the archived row above (same
# column layout as "New") before typing anything into them.
$wsOld.Range("A207:F207").Copy()
$wsNew.Range("A2:F10").PasteSpecial(-4122)   # xlPasteFormats
$wsNew.Application.CutCopyMode = $false

$wsNew.Range("B2:B10").NumberFormat = "@"
$wsNew.Range("E2:E10").NumberFormat = "@"

$row = 2
foreach ($r in $newRows) {
    $wsNew.Range("A" + $row).Value = $r.Link
    $wsNew.Range("B" + $row).Value = $r.Price
    $wsNew.Range("C" + $row).Value = $r.District
    $wsNew.Range("D" + $row).Value = $r.Area
    $wsNew.Range("E" + $row).Value = $r.Cadastre
    $wsNew.Range("F" + $row).Value = $r.Date
    $wsNew.Hyperlinks.Add($wsNew.Range("A" + $row), $r.Link)
    $row = $row + 1
}

# Same restoration as above: undo the built-in "Hyperlink" look that
# Hyperlinks.Add just stamped onto column A.
$wsOld.Range("A207").Copy()
$wsNew.Range("A2:A10").PasteSpecial(-4122)   # xlPasteFormats
$wsNew.Application.CutCopyMode = $false
